# Auto-generated edit script: rebuild Overview / zh-cn / de-de sheets
$wb = $excel.ActiveWorkbook

# ===== Overview sheet =====
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Hyperlinks.Delete()

$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"
$ws1.Range("D2").Value = "2016-03-25 00:28:01"
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/23397cf1dcb26096538c3a4049e9c41672faeb41/e2e/9bbe4cce-bbae-4080-bad5-65e2d9ab3115.md", "", "", "9bbe4cce-bbae-4080-bad5-65e2d9ab3115.md") | Out-Null

$ws1.Range("B3").Value = "Handed back: in sync with en-US"
$ws1.Range("C3").Value = "Handed back: in sync with en-US"
$ws1.Range("D3").Value = "2016-03-25 00:28:01"
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/23397cf1dcb26096538c3a4049e9c41672faeb41/e2e/d0ff0970-0a1d-4feb-b96c-3925dd614049.md", "", "", "d0ff0970-0a1d-4feb-b96c-3925dd614049.md") | Out-Null

$ws1.Range("B4").Value = "In Translation"
$ws1.Range("C4").Value = "In Translation"
$ws1.Range("D4").Value = "2016-03-25 00:25:44"
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/a65a288ad1b32b1b7b15436e6cc03a568c5a6953/e2e/d656325f-a872-42bc-9b12-be7455269f4e.md", "", "", "d656325f-a872-42bc-9b12-be7455269f4e.md") | Out-Null

$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"
$ws1.Range("D5").Value = "2016-03-25 00:28:01"
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/a451181bbd8179812b0969970f42aae01f4d4058/e2e/011fd86f-5c9c-45ef-ad39-97dcb22d32c3.md", "", "", "011fd86f-5c9c-45ef-ad39-97dcb22d32c3.md") | Out-Null

# ===== zh-cn sheet =====
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Hyperlinks.Delete()

$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/23397cf1dcb26096538c3a4049e9c41672faeb41/e2e/9bbe4cce-bbae-4080-bad5-65e2d9ab3115.md", "", "", "9bbe4cce-bbae-4080-bad5-65e2d9ab3115.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0902e3be3a691a857abe743a1427aec829d1248f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/9bbe4cce-bbae-4080-bad5-65e2d9ab3115.bf67642933cf71954002e7daae6e8c1b113c711e.zh-cn.xlf", "", "", "9bbe4cce-bbae-4080-bad5-65e2d9ab3115.bf67642933cf71954002e7daae6e8c1b113c711e.zh-cn.xlf") | Out-Null
$ws2.Range("E2").Value = "2016-03-25 00:27:57"
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/23397cf1dcb26096538c3a4049e9c41672faeb41/e2e/9bbe4cce-bbae-4080-bad5-65e2d9ab3115.md", "", "", "9bbe4cce-bbae-4080-bad5-65e2d9ab3115.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0902e3be3a691a857abe743a1427aec829d1248f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/9bbe4cce-bbae-4080-bad5-65e2d9ab3115.bf67642933cf71954002e7daae6e8c1b113c711e.zh-cn.xlf", "", "", "9bbe4cce-bbae-4080-bad5-65e2d9ab3115.bf67642933cf71954002e7daae6e8c1b113c711e.zh-cn.xlf") | Out-Null
$ws2.Range("H2").Value = "2016-03-25 00:28:23"
$ws2.Range("J2").Value = "Include"

$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/23397cf1dcb26096538c3a4049e9c41672faeb41/e2e/d0ff0970-0a1d-4feb-b96c-3925dd614049.md", "", "", "d0ff0970-0a1d-4feb-b96c-3925dd614049.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0902e3be3a691a857abe743a1427aec829d1248f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/d0ff0970-0a1d-4feb-b96c-3925dd614049.8c89534bfaae6eab7590d6d47039bead4084da69.zh-cn.xlf", "", "", "d0ff0970-0a1d-4feb-b96c-3925dd614049.8c89534bfaae6eab7590d6d47039bead4084da69.zh-cn.xlf") | Out-Null
$ws2.Range("E3").Value = "2016-03-25 00:27:57"
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/23397cf1dcb26096538c3a4049e9c41672faeb41/e2e/d0ff0970-0a1d-4feb-b96c-3925dd614049.md", "", "", "d0ff0970-0a1d-4feb-b96c-3925dd614049.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0902e3be3a691a857abe743a1427aec829d1248f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/d0ff0970-0a1d-4feb-b96c-3925dd614049.8c89534bfaae6eab7590d6d47039bead4084da69.zh-cn.xlf", "", "", "d0ff0970-0a1d-4feb-b96c-3925dd614049.8c89534bfaae6eab7590d6d47039bead4084da69.zh-cn.xlf") | Out-Null
$ws2.Range("H3").Value = "2016-03-25 00:28:23"
$ws2.Range("J3").Value = "Include"

$ws2.Range("B4").Value = ".md"
$ws2.Range("C4").Value = "In Translation"
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/a65a288ad1b32b1b7b15436e6cc03a568c5a6953/e2e/d656325f-a872-42bc-9b12-be7455269f4e.md", "", "", "d656325f-a872-42bc-9b12-be7455269f4e.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ed0790653e430b3f2a3619d77ee5436c42fc5ec8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d656325f-a872-42bc-9b12-be7455269f4e.55e410e4e13847dc2411b2130a62324de8fe6920.zh-cn.xlf", "", "", "d656325f-a872-42bc-9b12-be7455269f4e.55e410e4e13847dc2411b2130a62324de8fe6920.zh-cn.xlf") | Out-Null
$ws2.Range("E4").Value = "2016-03-25 00:25:40"
$ws2.Range("H4").Value = "0001-01-01 00:00:00"
$ws2.Range("J4").Value = "Include"

$ws2.Range("B5").Value = ".md"
$ws2.Range("C5").Value = "Ready for handoff"
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/a451181bbd8179812b0969970f42aae01f4d4058/e2e/011fd86f-5c9c-45ef-ad39-97dcb22d32c3.md", "", "", "011fd86f-5c9c-45ef-ad39-97dcb22d32c3.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0902e3be3a691a857abe743a1427aec829d1248f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/011fd86f-5c9c-45ef-ad39-97dcb22d32c3.9298244a657808e7b8bc3edcd758ef9e6446dd36.zh-cn.xlf", "", "", "011fd86f-5c9c-45ef-ad39-97dcb22d32c3.9298244a657808e7b8bc3edcd758ef9e6446dd36.zh-cn.xlf") | Out-Null
$ws2.Range("E5").Value = "2016-03-25 00:27:57"
$ws2.Range("H5").Value = "0001-01-01 00:00:00"
$ws2.Range("J5").Value = "Include"

# ===== de-de sheet =====
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Hyperlinks.Delete()

$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/23397cf1dcb26096538c3a4049e9c41672faeb41/e2e/9bbe4cce-bbae-4080-bad5-65e2d9ab3115.md", "", "", "9bbe4cce-bbae-4080-bad5-65e2d9ab3115.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9b5c3d88b691be2fa529f55021af2cd52a7e5b1b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/9bbe4cce-bbae-4080-bad5-65e2d9ab3115.bf67642933cf71954002e7daae6e8c1b113c711e.de-de.xlf", "", "", "9bbe4cce-bbae-4080-bad5-65e2d9ab3115.bf67642933cf71954002e7daae6e8c1b113c711e.de-de.xlf") | Out-Null
$ws3.Range("E2").Value = "2016-03-25 00:28:01"
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/23397cf1dcb26096538c3a4049e9c41672faeb41/e2e/9bbe4cce-bbae-4080-bad5-65e2d9ab3115.md", "", "", "9bbe4cce-bbae-4080-bad5-65e2d9ab3115.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9b5c3d88b691be2fa529f55021af2cd52a7e5b1b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/9bbe4cce-bbae-4080-bad5-65e2d9ab3115.bf67642933cf71954002e7daae6e8c1b113c711e.de-de.xlf", "", "", "9bbe4cce-bbae-4080-bad5-65e2d9ab3115.bf67642933cf71954002e7daae6e8c1b113c711e.de-de.xlf") | Out-Null
$ws3.Range("H2").Value = "2016-03-25 00:28:30"
$ws3.Range("J2").Value = "Include"

$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/23397cf1dcb26096538c3a4049e9c41672faeb41/e2e/d0ff0970-0a1d-4feb-b96c-3925dd614049.md", "", "", "d0ff0970-0a1d-4feb-b96c-3925dd614049.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9b5c3d88b691be2fa529f55021af2cd52a7e5b1b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/d0ff0970-0a1d-4feb-b96c-3925dd614049.8c89534bfaae6eab7590d6d47039bead4084da69.de-de.xlf", "", "", "d0ff0970-0a1d-4feb-b96c-3925dd614049.8c89534bfaae6eab7590d6d47039bead4084da69.de-de.xlf") | Out-Null
$ws3.Range("E3").Value = "2016-03-25 00:28:01"
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/23397cf1dcb26096538c3a4049e9c41672faeb41/e2e/d0ff0970-0a1d-4feb-b96c-3925dd614049.md", "", "", "d0ff0970-0a1d-4feb-b96c-3925dd614049.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9b5c3d88b691be2fa529f55021af2cd52a7e5b1b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/d0ff0970-0a1d-4feb-b96c-3925dd614049.8c89534bfaae6eab7590d6d47039bead4084da69.de-de.xlf", "", "", "d0ff0970-0a1d-4feb-b96c-3925dd614049.8c89534bfaae6eab7590d6d47039bead4084da69.de-de.xlf") | Out-Null
$ws3.Range("H3").Value = "2016-03-25 00:28:30"
$ws3.Range("J3").Value = "Include"

$ws3.Range("B4").Value = ".md"
$ws3.Range("C4").Value = "In Translation"
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/a65a288ad1b32b1b7b15436e6cc03a568c5a6953/e2e/d656325f-a872-42bc-9b12-be7455269f4e.md", "", "", "d656325f-a872-42bc-9b12-be7455269f4e.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/44e2b92f41e2d4d7bc19aa7e1816c673798cb51a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d656325f-a872-42bc-9b12-be7455269f4e.55e410e4e13847dc2411b2130a62324de8fe6920.de-de.xlf", "", "", "d656325f-a872-42bc-9b12-be7455269f4e.55e410e4e13847dc2411b2130a62324de8fe6920.de-de.xlf") | Out-Null
$ws3.Range("E4").Value = "2016-03-25 00:25:44"
$ws3.Range("H4").Value = "0001-01-01 00:00:00"
$ws3.Range("J4").Value = "Include"

$ws3.Range("B5").Value = ".md"
$ws3.Range("C5").Value = "Ready for handoff"
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/a451181bbd8179812b0969970f42aae01f4d4058/e2e/011fd86f-5c9c-45ef-ad39-97dcb22d32c3.md", "", "", "011fd86f-5c9c-45ef-ad39-97dcb22d32c3.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9b5c3d88b691be2fa529f55021af2cd52a7e5b1b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/011fd86f-5c9c-45ef-ad39-97dcb22d32c3.9298244a657808e7b8bc3edcd758ef9e6446dd36.de-de.xlf", "", "", "011fd86f-5c9c-45ef-ad39-97dcb22d32c3.9298244a657808e7b8bc3edcd758ef9e6446dd36.de-de.xlf") | Out-Null
$ws3.Range("E5").Value = "2016-03-25 00:28:01"
$ws3.Range("H5").Value = "0001-01-01 00:00:00"
$ws3.Range("J5").Value = "Include"

